$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 2; $r -le 19; $r++) {
    $cell = $ws.Range("C$r")
    if ($cell.Value2 -eq 45175) {
        $cell.Value2 = 45183
    }
}
